$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2744.35
$ws.Range("I40").Value = 2682.818
$ws.Range("J40").Value = 2819.5557
$ws.Range("K40").Value = 2682.818
$ws.Range("L40").Value = 2819.5557
$ws.Range("M40").Value = -2507.818
$ws.Range("N40").Value = -3169.5557

$ws.Range("H107").Value = 388.80768
$ws.Range("I107").Value = 216.80952
$ws.Range("J107").Value = 1111.2
$ws.Range("K107").Value = 216.80952
$ws.Range("L107").Value = 1111.2
$ws.Range("M107").Value = 1703.19048
$ws.Range("N107").Value = -4951.2

$ws.Range("H125").Value = 3955185
$ws.Range("J125").Value = 5882777.5
$ws.Range("L125").Value = 52944997.5
$ws.Range("N125").Value = -52949917.5

$ws.Range("H135").Value = 1204220.8
$ws.Range("I135").Value = 1009.4762
$ws.Range("J135").Value = 4011713.8
$ws.Range("K135").Value = 9085.2858
$ws.Range("L135").Value = 36105424.2
$ws.Range("M135").Value = -6550.2858
$ws.Range("N135").Value = -36110494.2

$ws.Range("H138").Value = 2773.05
$ws.Range("I138").Value = 1601.4642
$ws.Range("J138").Value = 4264.159
$ws.Range("K138").Value = 4804.392599999999
$ws.Range("L138").Value = 12792.477
$ws.Range("M138").Value = 335.6074000000008
$ws.Range("N138").Value = -23072.477

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2373.625
$ws.Range("I2").Value = 2427
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = 2427
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = -2314
$ws.Range("N2").Value = -2226

$ws.Range("H32").Value = 3773.11
$ws.Range("I32").Value = 3773.11
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 3773.11
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -3486.11
$ws.Range("N32").ClearContents()

$ws.Range("H45").Value = 2199.8333
$ws.Range("I45").Value = 1199.8572
$ws.Range("K45").Value = 1199.8572
$ws.Range("M45").Value = -822.8571999999999

$ws.Range("H74").Value = 1468.841
$ws.Range("I74").Value = 770.5
$ws.Range("J74").Value = 2965.2856
$ws.Range("K74").Value = 770.5
$ws.Range("L74").Value = 2965.2856
$ws.Range("M74").Value = 103.5
$ws.Range("N74").Value = -4713.2856

$ws.Range("H77").Value = 1468.841
$ws.Range("I77").Value = 770.5
$ws.Range("J77").Value = 2965.2856
$ws.Range("K77").Value = 3852.5
$ws.Range("L77").Value = 14826.428
$ws.Range("M77").Value = 515.5
$ws.Range("N77").Value = -23562.428

$ws.Range("H116").Value = 2373.625
$ws.Range("I116").Value = 2427
$ws.Range("J116").Value = 2000
$ws.Range("K116").Value = 2427
$ws.Range("L116").Value = 2000
$ws.Range("M116").Value = -133
$ws.Range("N116").Value = -6588

$ws.Range("H122").Value = 2170.4285
$ws.Range("I122").Value = 1430.2222
$ws.Range("J122").Value = 3502.8
$ws.Range("K122").Value = 4290.6666
$ws.Range("L122").Value = 10508.4
$ws.Range("M122").Value = -1840.6666
$ws.Range("N122").Value = -15408.4

$ws.Range("H132").Value = 1450.5223
$ws.Range("I132").Value = 1093.0613
$ws.Range("J132").Value = 2423.611
$ws.Range("K132").Value = 3279.1839
$ws.Range("L132").Value = 7270.833
$ws.Range("M132").Value = -749.1839
$ws.Range("N132").Value = -12330.833

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2373.625
$ws.Range("I3").Value = 2427
$ws.Range("J3").Value = 2000
$ws.Range("K3").Value = 2427
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = -2313
$ws.Range("N3").Value = -2228

$ws.Range("H94").Value = 784.46155
$ws.Range("I94").Value = 1105.4286
$ws.Range("J94").Value = 410
$ws.Range("K94").Value = 1105.4286
$ws.Range("L94").Value = 410
$ws.Range("M94").Value = -654.4286
$ws.Range("N94").Value = -1312

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1526.93
$ws.Range("I31").Value = 1151.5322
$ws.Range("J31").Value = 2139.4211
$ws.Range("K31").Value = 1151.5322
$ws.Range("L31").Value = 2139.4211
$ws.Range("M31").Value = -856.5322000000001
$ws.Range("N31").Value = -2729.4211

$ws.Range("H34").Value = 1526.93
$ws.Range("I34").Value = 1151.5322
$ws.Range("J34").Value = 2139.4211
$ws.Range("K34").Value = 1151.5322
$ws.Range("L34").Value = 2139.4211
$ws.Range("M34").Value = -949.5322000000001
$ws.Range("N34").Value = -2543.4211

$ws.Range("H62").Value = 4925
$ws.Range("I62").Value = 2433.3333
$ws.Range("J62").Value = 6420
$ws.Range("K62").Value = 2433.3333
$ws.Range("L62").Value = 6420
$ws.Range("M62").Value = -1809.3333
$ws.Range("N62").Value = -7668

$ws.Range("H65").Value = 4925
$ws.Range("I65").Value = 2433.3333
$ws.Range("J65").Value = 6420
$ws.Range("K65").Value = 12166.6665
$ws.Range("L65").Value = 32100
$ws.Range("M65").Value = -9046.666499999999
$ws.Range("N65").Value = -38340

$ws.Range("H86").Value = 27782400
$ws.Range("I86").Value = 45458544
$ws.Range("K86").Value = 45458544
$ws.Range("M86").Value = -45457421

$ws.Range("H89").Value = 27782400
$ws.Range("I89").Value = 45458544
$ws.Range("K89").Value = 227292720
$ws.Range("M89").Value = -227287104

$ws.Range("H99").Value = 1530983
$ws.Range("I99").Value = 2135449.5
$ws.Range("J99").Value = 19816.666
$ws.Range("K99").Value = 2135449.5
$ws.Range("L99").Value = 19816.666
$ws.Range("M99").Value = -2133951.5
$ws.Range("N99").Value = -22812.666

$ws.Range("H126").Value = 1530983
$ws.Range("I126").Value = 2135449.5
$ws.Range("J126").Value = 19816.666
$ws.Range("K126").Value = 6406348.5
$ws.Range("L126").Value = 59449.99800000001
$ws.Range("M126").Value = -6403878.5
$ws.Range("N126").Value = -64389.99800000001

$ws.Range("H132").Value = 759678.1
$ws.Range("I132").Value = 1191680.5
$ws.Range("J132").Value = 3673.9167
$ws.Range("K132").Value = 3575041.5
$ws.Range("L132").Value = 11021.7501
$ws.Range("M132").Value = -3572511.5
$ws.Range("N132").Value = -16081.7501

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1533.7941
$ws.Range("I129").Value = 750
$ws.Range("J129").Value = 1961.3182
$ws.Range("K129").Value = 2250
$ws.Range("L129").Value = 5883.9546
$ws.Range("M129").Value = 2750
$ws.Range("N129").Value = -15883.9546

$ws.Range("H132").Value = 2268.0588
$ws.Range("I132").Value = 2402.6667
$ws.Range("J132").Value = 2239.2144
$ws.Range("K132").Value = 21624.0003
$ws.Range("L132").Value = 20152.9296
$ws.Range("M132").Value = -19094.0003
$ws.Range("N132").Value = -25212.9296

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1311.9
$ws.Range("I97").Value = 1372.8235
$ws.Range("J97").Value = 966.6667
$ws.Range("K97").Value = 1372.8235
$ws.Range("L97").Value = 966.6667
$ws.Range("M97").Value = -876.8235
$ws.Range("N97").Value = -1958.6667

$ws.Range("H107").Value = 1496
$ws.Range("I107").Value = 2717.5
$ws.Range("J107").Value = 518.8
$ws.Range("K107").Value = 2717.5
$ws.Range("L107").Value = 518.8
$ws.Range("M107").Value = -797.5
$ws.Range("N107").Value = -4358.8

$ws.Range("H113").Value = 1911
$ws.Range("I113").Value = 1680
$ws.Range("J113").Value = 2296
$ws.Range("K113").Value = 1680
$ws.Range("L113").Value = 2296
$ws.Range("M113").Value = 490
$ws.Range("N113").Value = -6636

$ws.Range("H132").Value = 1672.7646
$ws.Range("I132").Value = 1332.5897
$ws.Range("J132").Value = 2778.3333
$ws.Range("K132").Value = 3997.7691
$ws.Range("L132").Value = 8334.999899999999
$ws.Range("M132").Value = -1467.7691
$ws.Range("N132").Value = -13394.9999

$ws.Range("H135").Value = 41175.8
$ws.Range("J135").Value = 41175.8
$ws.Range("L135").Value = 41175.8
$ws.Range("N135").Value = -51315.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 80042.30499999999
$ws.Range("I7").Value = 102780
$ws.Range("J7").Value = 4250
$ws.Range("K7").Value = 102780
$ws.Range("L7").Value = 4250
$ws.Range("M7").Value = -102668
$ws.Range("N7").Value = -4474

$ws.Range("H68").Value = 1995.6
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 1995.6
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 1995.6
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -3493.6

$ws.Range("H71").Value = 1995.6
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 1995.6
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 9978
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -17466

$ws.Range("H122").Value = 37038600
$ws.Range("I122").Value = 55556560
$ws.Range("J122").Value = 2680
$ws.Range("K122").Value = 166669680
$ws.Range("L122").Value = 8040
$ws.Range("M122").Value = -166667230
$ws.Range("N122").Value = -12940

$ws.Range("H126").Value = 80042.30499999999
$ws.Range("I126").Value = 102780
$ws.Range("J126").Value = 4250
$ws.Range("K126").Value = 308340
$ws.Range("L126").Value = 12750
$ws.Range("M126").Value = -305870
$ws.Range("N126").Value = -17690

$ws.Range("H132").Value = 6450.7793
$ws.Range("I132").Value = 6606.0654
$ws.Range("J132").Value = 6126.091
$ws.Range("K132").Value = 19818.1962
$ws.Range("L132").Value = 18378.273
$ws.Range("M132").Value = -17288.1962
$ws.Range("N132").Value = -23438.273

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1000
$ws.Range("I96").Value = 1000
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 1000
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = 373
$ws.Range("N96").ClearContents()

$ws.Range("H122").Value = 49486
$ws.Range("I122").Value = 64175.375
$ws.Range("J122").Value = 2480
$ws.Range("K122").Value = 192526.125
$ws.Range("L122").Value = 7440
$ws.Range("M122").Value = -190076.125
$ws.Range("N122").Value = -12340

$ws.Range("H132").Value = 1644.75
$ws.Range("I132").Value = 1122.8529
$ws.Range("J132").Value = 2327.2307
$ws.Range("K132").Value = 3368.5587
$ws.Range("L132").Value = 6981.6921
$ws.Range("M132").Value = -838.5587000000005
$ws.Range("N132").Value = -12041.6921
